# Natmi following Dr Hou advice
# Rewrites the LR-pairs result rows for Ntf3-Ntrk2: existing rows 2-7 get
# updated statistics (ligand-expressing-cell counts, weights, specificities)
# and three new rows (8-10) are added for the "sCs" sending cluster, so the
# used range grows from A1:T7 to A1:T10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Each inner array holds columns A..T for one data row (row numbers 2..10).
# Columns A-D are strings (Sending cluster, Ligand symbol, Receptor symbol,
# Target cluster); columns E-T are the numeric NATMI statistics.
$data = @(
  # row 2: ECs -> ECs
  @("ECs","Ntf3","Ntrk2","ECs",
    3,1,3.998902333333333,11.996707,0.4953291676276626,0.4953291676276626,
    3,1,0.5901906666666666,1.770572,0.05042357210835777,0.05042357210835777,
    2.360114834044889,21.241033506404,0.02497626600124628,0.02497626600124628),

  # row 3: ECs -> FAPs
  @("ECs","Ntf3","Ntrk2","FAPs",
    3,1,3.998902333333333,11.996707,0.4953291676276626,0.4953291676276626,
    3,1,10.74310933333333,32.229328,0.9178490591808264,0.9178490591808263,
    42.96064498032177,386.6458048228959,0.454637410491872,0.4546374104918719),

  # row 4: ECs -> sCs
  @("ECs","Ntf3","Ntrk2","sCs",
    3,1,3.998902333333333,11.996707,0.4953291676276626,0.4953291676276626,
    3,1,0.371358,1.114074,0.03172736871081582,0.03172736871081581,
    1.485024372702,13.365219354318,0.01571549113454435,0.01571549113454435),

  # row 5: FAPs -> ECs
  @("FAPs","Ntf3","Ntrk2","ECs",
    3,1,3.99311,11.97933,0.4946116928284643,0.4946116928284644,
    3,1,0.5901906666666666,1.770572,0.05042357210835777,0.05042357210835777,
    2.356696252973333,21.21026627676,0.02494008835897298,0.02494008835897298),

  # row 6: FAPs -> FAPs
  @("FAPs","Ntf3","Ntrk2","FAPs",
    3,1,3.99311,11.97933,0.4946116928284643,0.4946116928284644,
    3,1,10.74310933333333,32.229328,0.9178490591808264,0.9178490591808263,
    42.89841731002666,386.08575579024,0.4539788769224419,0.4539788769224419),

  # row 7: FAPs -> sCs
  @("FAPs","Ntf3","Ntrk2","sCs",
    3,1,3.99311,11.97933,0.4946116928284643,0.4946116928284644,
    3,1,0.371358,1.114074,0.03172736871081582,0.03172736871081581,
    1.48287334338,13.34586009042,0.01569272754704946,0.01569272754704946),

  # row 8 (new): sCs -> ECs
  @("sCs","Ntf3","Ntrk2","ECs",
    1,0.3333333333333333,0.08120966666666667,0.243629,0.01005913954387315,0.01005913954387315,
    3,1,0.5901906666666666,1.770572,0.05042357210835777,0.05042357210835777,
    0.04792918730977778,0.431362685788,0.0005072177481385209,0.0005072177481385209),

  # row 9 (new): sCs -> FAPs
  @("sCs","Ntf3","Ntrk2","FAPs",
    1,0.3333333333333333,0.08120966666666667,0.243629,0.01005913954387315,0.01005913954387315,
    3,1,10.74310933333333,32.229328,0.9178490591808264,0.9178490591808263,
    0.8724443279235554,7.851998951311999,0.009232771766512618,0.009232771766512616),

  # row 10 (new): sCs -> sCs
  @("sCs","Ntf3","Ntrk2","sCs",
    1,0.3333333333333333,0.08120966666666667,0.243629,0.01005913954387315,0.01005913954387315,
    3,1,0.371358,1.114074,0.03172736871081582,0.03172736871081581,
    0.030157859394,0.271420734546,0.0003191500292220111,0.000319150029222011)
)

for ($r = 0; $r -lt $data.Length; $r++) {
  $rowNum = $r + 2
  $rowData = $data[$r]
  for ($c = 0; $c -lt $cols.Length; $c++) {
    $ws.Range($cols[$c] + $rowNum).Value = $rowData[$c]
  }
}
